$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row 43 (pushes nothing below it since it's the last row),
# carrying over row 42's formatting - this keeps the date-styled A column
# using the existing style (s="1") instead of minting a new cellXf.
$ws.Rows(43).Insert()

$ws.Range("A43").Value = 45809
$ws.Range("B43").Value = -0.093
$ws.Range("C43").Value = 0.623
$ws.Range("D43").Value = -0.178
$ws.Range("E43").Value = 0.342
$ws.Range("F43").Value = 0.904
